# Insert a new weekly price-observation row for Piña (Caramelo / Segunda)
# at sheet row 338, pushing the existing rows 338:366 down to 339:367.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("338:338").Insert()

$ws.Range("A338").Value = 5
$ws.Range("B338").Value = "Macroferia Regional de Talca"
$ws.Range("C338").Value = "Maule"
$ws.Range("D338").Value = 45013
$ws.Range("E338").Value = 7
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100108
$ws.Range("H338").Value = "Tropicales y subtropicales"
$ws.Range("I338").Value = 100108005
$ws.Range("J338").Value = "Piña"
$ws.Range("K338").Value = "Caramelo"
$ws.Range("L338").Value = "Segunda"
$ws.Range("M338").Value = 300
$ws.Range("N338").Value = 20000
$ws.Range("O338").Value = 20000
$ws.Range("P338").Value = 20000
$ws.Range("Q338").Value = "$/caja 14 unidades"
$ws.Range("R338").Value = "Ecuador"
$ws.Range("S338").Value = 1429
$ws.Range("T338").Value = 14
